$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row for 6th April 2020
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A24").Value = 43927
$ws.Range("B24").Value = 16
$ws.Range("C24").Value = 4277
$ws.Range("F24").Value = 158
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = "Uhuru"

# Update view: scroll so column D is leftmost, select L24
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("L24").Select()
